$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Change 1: paragraph w14:paraId="0000000D" (the <personal-statement>
# body text paragraph).
#   - pPr gains <w:ind w:left="425.19685039370086"
#     w:hanging="850.3937007874017"/>
#   - pPr/rPr gains <w:color w:val="17365d"/>
#   - the single text run is split into a leading run of 13 spaces plus
#     the (edited) text run
#   - the trailing empty run is dropped
#   - wording tweaks: "during  any practical work" -> "during my work",
#     "is an opportunity to get" -> "is to have an opportunity to gain"
# ----------------------------------------------------------------------

$newPara1 = '<w:p w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="0000000D"><w:pPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="425.19685039370086" w:hanging="850.3937007874017"/><w:rPr><w:rFonts w:ascii="Yu Gothic Medium" w:cs="Yu Gothic Medium" w:eastAsia="Yu Gothic Medium" w:hAnsi="Yu Gothic Medium"/><w:color w:val="17365d"/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Yu Gothic Medium" w:cs="Yu Gothic Medium" w:eastAsia="Yu Gothic Medium" w:hAnsi="Yu Gothic Medium"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">             </w:t></w:r><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Yu Gothic Medium" w:cs="Yu Gothic Medium" w:eastAsia="Yu Gothic Medium" w:hAnsi="Yu Gothic Medium"/><w:color w:val="17365d"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">I am a Front-end developer. I am ready to improve  my knowledge and increase it during my work. I am considering reaching success in the field of IT so I have chosen just this branch as my speciality. Currently, the priority for me is to have an opportunity to gain new knowledge working with a qualified and creative team.</w:t></w:r></w:p>'

$found1 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^\s*I am a Front-end developer") {
        [void]$p.Range.InsertXML($newPara1)
        $found1 = $true
        break
    }
}
if (-not $found1) {
    throw "Change 1: target paragraph not found"
}

# ----------------------------------------------------------------------
# Change 2: paragraph w14:paraId="0000000F" (the blank line right after
# the closing </personal-statement> paragraph) gains
# <w:color w:val="ffffff"/> inside its pPr/rPr.
# ----------------------------------------------------------------------

$newPara2 = '<w:p w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="0000000F"><w:pPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="425.19685039370086" w:hanging="850.3937007874017"/><w:rPr><w:rFonts w:ascii="Yu Gothic Medium" w:cs="Yu Gothic Medium" w:eastAsia="Yu Gothic Medium" w:hAnsi="Yu Gothic Medium"/><w:color w:val="ffffff"/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'

$found2 = $false
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $prevText = $d.Paragraphs.Item($i - 1).Range.Text
    if ($t -match "^\s*$" -and $prevText -match "</personal-statement>") {
        # the blank paragraph immediately following the closing
        # </personal-statement> tag
        [void]$p.Range.InsertXML($newPara2)
        $found2 = $true
        break
    }
}
if (-not $found2) {
    throw "Change 2: target paragraph not found"
}
